# New weekly "Betarraga" price record (Primera/Segunda) is inserted at the
# top of the data block (row 90), pushing the existing rows 90-178 down to
# 92-180. Excel's default "insert, shift down" behaviour copies the row
# above's formatting (the date-format style "s=2" on column D), so nothing
# further needs to be done for styling - only the 18 cell values for the
# two new rows need to be populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the existing row 90, shifting everything
# from row 90 down to row 92 (and the old last row, 178, down to 180).
$ws.Rows("90:91").Insert()

# Row 90 - "Primera" quality for the new date.
$ws.Range("A90").Value = 1
$ws.Range("B90").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C90").Value = "Arica y Parinacota"
$ws.Range("D90").Value = 44484
$ws.Range("E90").Value = 15
$ws.Range("F90").Value = 100114014
$ws.Range("G90").Value = "Betarraga"
$ws.Range("H90").Value = "Sin especificar"
$ws.Range("I90").Value = "Primera"
$ws.Range("J90").Value = 1200
$ws.Range("K90").Value = 350
$ws.Range("L90").Value = 400
$ws.Range("M90").Value = 375
$ws.Range("N90").Value = "$/paquete 4 unidades"
$ws.Range("O90").Value = "Región de Arica y Parinacota"
$ws.Range("P90").Value = 94
$ws.Range("Q90").Value = 4
$ws.Range("R90").Value = "Hortaliza"

# Row 91 - "Segunda" quality for the new date.
$ws.Range("A91").Value = 1
$ws.Range("B91").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C91").Value = "Arica y Parinacota"
$ws.Range("D91").Value = 44484
$ws.Range("E91").Value = 15
$ws.Range("F91").Value = 100114014
$ws.Range("G91").Value = "Betarraga"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Segunda"
$ws.Range("J91").Value = 1000
$ws.Range("K91").Value = 350
$ws.Range("L91").Value = 400
$ws.Range("M91").Value = 375
$ws.Range("N91").Value = "$/paquete 5 unidades"
$ws.Range("O91").Value = "Región de Arica y Parinacota"
$ws.Range("P91").Value = 75
$ws.Range("Q91").Value = 5
$ws.Range("R91").Value = "Hortaliza"
